$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "keyword" columns (H and J) so each bare noun value is now
# stored as a Python-list-literal-looking string "['Word','noun']" instead of
# just "Word". This reflects the Name class now tagging each keyword with its
# part of speech. ---

$pairs = @(
    @{ Row = 2;  H = 'Trees';  J = 'Great'   },
    @{ Row = 3;  H = 'Lord';   J = 'July'    },
    @{ Row = 4;  H = 'Menu';   J = 'Name'    },
    @{ Row = 5;  H = 'Life';   J = 'Garden'  },
    @{ Row = 6;  H = 'Lit';    J = 'Kingdom' },
    @{ Row = 7;  H = 'Plant';  J = 'Year'    },
    @{ Row = 8;  H = 'May';    J = 'World'   },
    @{ Row = 9;  H = 'Trees';  J = 'Don'     },
    @{ Row = 10; H = 'Award';  J = 'Plate'   },
    @{ Row = 11; H = 'Fossil'; J = 'Jung'    }
)

foreach ($pair in $pairs) {
    $r = $pair.Row
    $ws.Cells.Item($r, 8).Value  = "['" + $pair.H + "','noun']"
    $ws.Cells.Item($r, 10).Value = "['" + $pair.J + "','noun']"
}

# --- Column width changes ---
# (target stored widths are 16.6640625 / 14.33203125 / 16 "character" widths;
# ColumnWidth is quantized to the nearest 1/6 of a character by the host, so
# these are the closest settable values that reproduce that width.)
$ws.Columns.Item(2).ColumnWidth  = 15.833333333333334
$ws.Columns.Item(8).ColumnWidth  = 13.5
$ws.Columns.Item(10).ColumnWidth = 15.166666666666666

# --- Selection change ---
$ws.Range("H14").Select()
